# ============================================================================
# Adds a new "2022-Q1" fund-holdings sheet (positioned right before the
# "总计" summary sheet) and records its aggregate stats as a new top data
# row in "总计" (shifting the existing rows down).
# ============================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Locate the existing "总计" (summary) sheet and the last quarterly sheet
#    ("2021-Q4") that we use as a formatting template for the new sheet.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$templateSheet = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------------
# 2. Insert the new worksheet immediately before "总计" and rename it.
#    NOTE: after the insert, the position that $totalSheet was anchored to
#    now resolves to the newly-inserted sheet, so "总计" must be re-fetched
#    by name afterwards to get a reference to the original sheet again.
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"
$totalSheet = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------------
# 3. Fund holdings data for 2022-Q1 (index, code, name, size, stock
#    position, position ratio, held market value (100M yuan), position rank)
# ---------------------------------------------------------------------------
$fundData = @(
    @(0, "090018", "大成新锐产业混合", "125.72", "89.64", "5.31", "6.6757", 4),
    @(1, "001300", "大成睿景灵活配置混合A", "67.47", "89.89", "5.39", "3.6366", 4),
    @(2, "013435", "大成景气精选六个月持有混合A", "46.91", "87.13", "4.79", "2.2470", 5),
    @(3, "001301", "大成睿景灵活配置混合C", "27.47", "89.89", "5.39", "1.4806", 4),
    @(4, "012519", "大成核心趋势混合型证券投资基金A", "16.82", "87.87", "5.88", "0.9890", 4),
    @(5, "002258", "大成国企改革灵活配置混合", "17.17", "89.07", "5.57", "0.9564", 3),
    @(6, "008934", "大成科技消费股票A", "13.43", "83.81", "5.83", "0.7830", 4),
    @(7, "160918", "大成中小盘混合(LOF)A", "8.66", "71.76", "5.78", "0.5005", 1),
    @(8, "010826", "大成产业趋势混合A", "8.93", "91.32", "5.53", "0.4938", 5),
    @(9, "012184", "大成创新趋势混合型证券投资基金A", "8.05", "72.67", "5.79", "0.4661", 1),
    @(10, "013436", "大成景气精选六个月持有混合C", "7.85", "87.13", "4.79", "0.3760", 5),
    @(11, "008935", "大成科技消费股票C", "3.98", "83.81", "5.83", "0.2320", 4),
    @(12, "014185", "招商专精特新股票A", "8.37", "30.94", "2.42", "0.2026", 2),
    @(13, "012520", "大成核心趋势混合型证券投资基金D", "3.29", "87.87", "5.88", "0.1935", 4),
    @(14, "008274", "大成行业先锋混合A", "3.19", "73.98", "5.79", "0.1847", 1),
    @(15, "010827", "大成产业趋势混合C", "2.01", "91.32", "5.53", "0.1112", 5),
    @(16, "002945", "大成盛世精选灵活配置混合", "1.52", "70.48", "5.84", "0.0888", 1),
    @(17, "014186", "招商专精特新股票C", "3.46", "30.94", "2.42", "0.0837", 2),
    @(18, "008275", "大成行业先锋混合C", "0.45", "73.98", "5.79", "0.0261", 1),
    @(19, "006230", "鹏华研究驱动混合", "0.85", "86.99", "1.83", "0.0156", 10),
    @(20, "011254", "长江量化科技精选一个月滚动持有股票型发起式证券投资基金A", "0.64", "92.60", "1.83", "0.0117", 7),
    @(21, "080007", "长盛同鑫行业配置混合", "0.27", "87.31", "4.02", "0.0109", 4),
    @(22, "012185", "大成创新趋势混合型证券投资基金C", "0.14", "72.67", "5.79", "0.0081", 1),
    @(23, "620004", "金元顺安价值增长混合", "0.18", "87.28", "3.49", "0.0063", 5),
    @(24, "620002", "金元顺安成长动力混合", "0.16", "62.78", "3.84", "0.0061", 3),
    @(25, "011255", "长江量化科技精选一个月滚动持有股票型发起式证券投资基金C", "0.10", "92.60", "1.83", "0.0018", 7),
    @(26, "011159", "大成中小盘混合(LOF)C", "0.01", "71.76", "5.78", "0.0006", 1)
)
$rowCount = $fundData.Count
$lastRow = $rowCount + 1

# ---------------------------------------------------------------------------
# 4. Header row.
# ---------------------------------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------------
# 5. Data rows: column A (index) and H (rank) are numeric; B-G are text
#    (mirrors the source data, which stores these figures as plain strings).
# ---------------------------------------------------------------------------
$idxArr = New-Object 'object[,]' $rowCount,1
$textArr = New-Object 'object[,]' $rowCount,5
$rankArr = New-Object 'object[,]' $rowCount,1

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $fundData[$i]
    $idxArr[$i,0] = $r[0]
    $textArr[$i,0] = $r[1]
    $textArr[$i,1] = $r[2]
    $textArr[$i,2] = $r[3]
    $textArr[$i,3] = $r[4]
    $textArr[$i,4] = $r[5]
    $rankArr[$i,0] = $r[7]
}

$newSheet.Range("A2:A$lastRow").Value = $idxArr
$newSheet.Range("B2:F$lastRow").NumberFormat = "@"
$newSheet.Range("B2:F$lastRow").Value = $textArr
$newSheet.Range("G2:G$lastRow").NumberFormat = "@"
for ($i = 0; $i -lt $rowCount; $i++) {
    $newSheet.Cells.Item($i + 2, 7).Value = $fundData[$i][6]
}
$newSheet.Range("H2:H$lastRow").Value = $rankArr

# ---------------------------------------------------------------------------
# 6. Match the formatting of the other quarterly sheets: bold/centered/
#    bordered header row and index column (style used throughout the rest
#    of the workbook).
# ---------------------------------------------------------------------------
$templateSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$templateSheet.Range("A2:A$lastRow").Copy()
$newSheet.Range("A2:A$lastRow").PasteSpecial(-4122)

$newSheet.Range("A1").Select()

# ---------------------------------------------------------------------------
# 7. Insert a new top data row into "总计" for the 2022-Q1 aggregate and
#    re-point the index column values (0-based) for the whole table.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("A2:D2").Value = $totalSheet.Range("A3:D3").Value
$totalSheet.Range("A2:D2").NumberFormat = $totalSheet.Range("A3:D3").NumberFormat
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A2:D2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 27
$totalSheet.Range("D2").Value = 19.79

$totalSheet.Range("A1").Select()
